$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.432.62'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '1.798.58'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '224.18'
$ws.Range("E5").Value = '  -1.79%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.602'
$ws.Range("E6").Value = '  -1.81%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.27'
$ws.Range("E8").Value = '  +6.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.287'
$ws.Range("E9").Value = '  -5.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0666'
$ws.Range("E10").Value = '  -3.98%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0984'
$ws.Range("E11").Value = '  +1.97%  '
$ws.Range("D12").Value = '2.058.83'
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.88'
$ws.Range("E13").Value = '  -5.69%  '
$ws.Range("D14").Value = '1.800.25'
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("D15").Value = '34.409.53'
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("E16").Value = '  -4.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.35'
$ws.Range("E17").Value = '  -3.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.89'
$ws.Range("E18").Value = '  -2.89%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '238.35'
$ws.Range("E19").Value = '  -2.81%  '
$ws.Range("D20").Value = '0.0₃0764'
$ws.Range("E20").Value = '  -3.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.08'
$ws.Range("E21").Value = '  -4.70%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.06'
$ws.Range("E23").Value = '  -3.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.15'
$ws.Range("E24").Value = '  -2.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.78'
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.53'
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.66'
$ws.Range("E27").Value = '  -4.20%  '
$ws.Range("E28").Value = '  -2.26%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("E31").Value = '  -3.19%  '
$ws.Range("E32").Value = '  -3.51%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.82'
$ws.Range("E33").Value = '  -4.75%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.81'
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("E35").Value = '  -0.88%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.635'
$ws.Range("E36").Value = '  -5.33%  '
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.299.26'
$ws.Range("E37").Value = '  -6.81%  '
$ws.Range("E38").Value = '  -2.73%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.29'
$ws.Range("E39").Value = '  -6.57%  '
$ws.Range("B40").Value = 'HuobiToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.44'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("E41").Value = '  +1.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.33'
$ws.Range("E42").Value = '  -1.44%  '
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.940'
$ws.Range("E44").Value = '  -2.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.90'
$ws.Range("E45").Value = '  +3.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0515'
$ws.Range("E46").Value = '  +3.71%  '
$ws.Range("D47").Value = '1.959.61'
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.70'
$ws.Range("E48").Value = '  -5.63%  '
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '101.42'
$ws.Range("E50").Value = '  -2.76%  '
$ws.Range("E51").Value = '  -0.81%  '
